# Update leveling profit data across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled market-price refresh run.
$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 1753.2059
$ws.Cells.Item(98, 9).Value = 1741.8966
$ws.Cells.Item(98, 10).Value = 1818.8
$ws.Cells.Item(98, 11).Value = 1741.8966
$ws.Cells.Item(98, 12).Value = 1818.8
$ws.Cells.Item(98, 13).Value = -243.8966
$ws.Cells.Item(98, 14).Value = -4814.8

# ALC row 122
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(122, 8).Value = 1753.2059
$ws.Cells.Item(122, 9).Value = 1741.8966
$ws.Cells.Item(122, 10).Value = 1818.8
$ws.Cells.Item(122, 11).Value = 5225.6898
$ws.Cells.Item(122, 12).Value = 5456.4
$ws.Cells.Item(122, 13).Value = -2775.6898
$ws.Cells.Item(122, 14).Value = -10356.4

# ALC row 132
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 5129315
$ws.Cells.Item(132, 10).Value = 1753.2858
$ws.Cells.Item(132, 12).Value = 5259.857400000001
$ws.Cells.Item(132, 14).Value = -10319.8574

# ALC row 135
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(135, 8).Value = 605.1739
$ws.Cells.Item(135, 9).Value = 553.6842
$ws.Cells.Item(135, 10).Value = 849.75
$ws.Cells.Item(135, 11).Value = 4983.1578
$ws.Cells.Item(135, 12).Value = 7647.75
$ws.Cells.Item(135, 13).Value = -2448.1578
$ws.Cells.Item(135, 14).Value = -12717.75

# ALC row 137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 1243.9143
$ws.Cells.Item(137, 9).Value = 1027
$ws.Cells.Item(137, 11).Value = 3081
$ws.Cells.Item(137, 13).Value = -531

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 1972.7606
$ws.Cells.Item(138, 10).Value = 2670.182
$ws.Cells.Item(138, 12).Value = 8010.545999999999
$ws.Cells.Item(138, 14).Value = -18290.546

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 3015.0366
$ws.Cells.Item(32, 9).Value = 2575.7368
$ws.Cells.Item(32, 11).Value = 2575.7368
$ws.Cells.Item(32, 13).Value = -2288.7368

# ARM row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 27779724
$ws.Cells.Item(61, 9).Value = 15626484
$ws.Cells.Item(61, 11).Value = 15626484
$ws.Cells.Item(61, 13).Value = -15626272

# ARM row 101
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(101, 8).Value = 87498.164
$ws.Cells.Item(101, 10).Value = 87498.164
$ws.Cells.Item(101, 12).Value = 87498.164
$ws.Cells.Item(101, 14).Value = -93988.164

# ARM row 114
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(114, 8).Value = 7950
$ws.Cells.Item(114, 10).Value = 7950
$ws.Cells.Item(114, 12).Value = 7950
$ws.Cells.Item(114, 14).Value = -16628

# ARM row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 27779724
$ws.Cells.Item(136, 9).Value = 15626484
$ws.Cells.Item(136, 11).Value = 46879452
$ws.Cells.Item(136, 13).Value = -46876902

# BSM row 92
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(92, 8).Value = 19998.5
$ws.Cells.Item(92, 10).Value = 19998.5
$ws.Cells.Item(92, 12).Value = 19998.5
$ws.Cells.Item(92, 14).Value = -24990.5

# BSM row 105
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(105, 8).Value = 2425.5356
$ws.Cells.Item(105, 9).Value = 2429.8333
$ws.Cells.Item(105, 11).Value = 2429.8333
$ws.Cells.Item(105, 13).Value = -682.8332999999998

# CRP row 31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 1787341.1
$ws.Cells.Item(31, 9).Value = 2646720
$ws.Cells.Item(31, 11).Value = 2646720
$ws.Cells.Item(31, 13).Value = -2646425

# CRP row 34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 1787341.1
$ws.Cells.Item(34, 9).Value = 2646720
$ws.Cells.Item(34, 11).Value = 2646720
$ws.Cells.Item(34, 13).Value = -2646518

# CRP row 74
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(74, 8).Value = 34998.668
$ws.Cells.Item(74, 10).Value = 34998.668
$ws.Cells.Item(74, 12).Value = 34998.668
$ws.Cells.Item(74, 14).Value = -36746.668

# CRP row 77
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(77, 8).Value = 34998.668
$ws.Cells.Item(77, 10).Value = 34998.668
$ws.Cells.Item(77, 12).Value = 104996.004
$ws.Cells.Item(77, 14).Value = -113732.004

# CRP row 95
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(95, 8).Value = 26182.2
$ws.Cells.Item(95, 10).Value = 26182.2
$ws.Cells.Item(95, 12).Value = 26182.2
$ws.Cells.Item(95, 14).Value = -31674.2

# CRP row 96
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(96, 8).Value = 17524.666
$ws.Cells.Item(96, 10).Value = 17524.666
$ws.Cells.Item(96, 12).Value = 17524.666
$ws.Cells.Item(96, 14).Value = -23016.666

# CRP row 132
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value = 1304.6613
$ws.Cells.Item(132, 9).Value = 1016.46
$ws.Cells.Item(132, 10).Value = 2505.5
$ws.Cells.Item(132, 11).Value = 3049.38
$ws.Cells.Item(132, 12).Value = 7516.5
$ws.Cells.Item(132, 13).Value = -519.3800000000001
$ws.Cells.Item(132, 14).Value = -12576.5

# CRP row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 1478.5272
$ws.Cells.Item(134, 9).Value = 1305.3334
$ws.Cells.Item(134, 11).Value = 3916.0002
$ws.Cells.Item(134, 13).Value = -1381.0002

# CUL row 122
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(122, 8).Value = 917.08826
$ws.Cells.Item(122, 10).Value = 1190.9
$ws.Cells.Item(122, 12).Value = 10718.1
$ws.Cells.Item(122, 14).Value = -15618.1

# CUL row 123
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(123, 8).Value = 15000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 15000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 45000
$ws.Cells.Item(123, 13).ClearContents()
$ws.Cells.Item(123, 14).Value = -49900

# GSM row 92
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(92, 8).Value = 20847.555
$ws.Cells.Item(92, 10).Value = 20847.555
$ws.Cells.Item(92, 12).Value = 20847.555
$ws.Cells.Item(92, 14).Value = -24591.555

# GSM row 102
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2066.0356
$ws.Cells.Item(102, 9).Value = 1950.8096
$ws.Cells.Item(102, 11).Value = 1950.8096
$ws.Cells.Item(102, 13).Value = -328.8096

# GSM row 104
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(104, 8).Value = 49996
$ws.Cells.Item(104, 10).Value = 49996
$ws.Cells.Item(104, 12).Value = 49996
$ws.Cells.Item(104, 14).Value = -56984

# GSM row 135
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(135, 8).Value = 60714.285
$ws.Cells.Item(135, 10).Value = 60714.285
$ws.Cells.Item(135, 12).Value = 60714.285
$ws.Cells.Item(135, 14).Value = -70854.285

# LTW row 68
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(68, 8).Value = 1554.9445
$ws.Cells.Item(68, 9).Value = 1414.6666
$ws.Cells.Item(68, 11).Value = 1414.6666
$ws.Cells.Item(68, 13).Value = -665.6666

# LTW row 71
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(71, 8).Value = 1554.9445
$ws.Cells.Item(71, 9).Value = 1414.6666
$ws.Cells.Item(71, 11).Value = 7073.333000000001
$ws.Cells.Item(71, 13).Value = -3329.333000000001

# LTW row 101
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(101, 8).Value = 14180.25
$ws.Cells.Item(101, 10).Value = 14180.25
$ws.Cells.Item(101, 12).Value = 14180.25
$ws.Cells.Item(101, 14).Value = -20670.25

# LTW row 106
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(106, 8).Value = 19999
$ws.Cells.Item(106, 10).Value = 19999
$ws.Cells.Item(106, 12).Value = 19999
$ws.Cells.Item(106, 14).Value = -22523

# LTW row 110
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(110, 8).Value = 6950
$ws.Cells.Item(110, 10).Value = 6950
$ws.Cells.Item(110, 12).Value = 6950
$ws.Cells.Item(110, 14).Value = -15130

# WVR row 104
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(104, 8).Value = 14925
$ws.Cells.Item(104, 10).Value = 14925
$ws.Cells.Item(104, 12).Value = 14925
$ws.Cells.Item(104, 14).Value = -21913

# WVR row 105
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(105, 8).Value = 49997.668
$ws.Cells.Item(105, 10).Value = 49997.668
$ws.Cells.Item(105, 12).Value = 49997.668
$ws.Cells.Item(105, 14).Value = -56985.668

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 1055.762
$ws.Cells.Item(132, 9).Value = 808.7907
$ws.Cells.Item(132, 11).Value = 2426.3721
$ws.Cells.Item(132, 13).Value = 103.6279

# WVR row 136
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 10685683
$ws.Cells.Item(136, 10).Value = 1748.6364
$ws.Cells.Item(136, 12).Value = 5245.9092
$ws.Cells.Item(136, 14).Value = -10345.9092
